# Update the "Correspond Handoff Datetime" (D4) and
# "Correspond Handback DateTime" (G4) timestamps on the language report
# sheets, as part of regenerating the handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-27 02:25:27"
$wsZhCn.Range("G4").Value = "2016-01-27 02:26:17"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-27 02:25:39"
$wsDeDe.Range("G4").Value = "2016-01-27 02:26:38"
